# This edit adds a new weekly price record for "Coliflor" (Cauliflower) at
# "Terminal La Palmera de La Serena" into the date-sorted data block that
# occupies rows 1096-1197 of Sheet1. The new record is inserted logically at
# the top of the block (row 1096), so every existing record in that block
# shifts down by one row, and a brand-new row 1198 is appended holding what
# used to be the last row (1197) of the block.
#
# Columns A,B,C,E,F,G,H,N,O,Q,R are constant across the whole block (same
# market/category/etc.), so only D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) actually need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 1096
$lastRow  = 1197
$newLastRow = 1198

$colD = 4   # Fecha
$colI = 9   # Calidad
$colJ = 10  # Volumen
$colK = 11  # Precio minimo
$colL = 12  # Precio maximo
$colM = 13  # Precio promedio ponderado
$colP = 16  # Precio $/Kg

# --- 1) snapshot the current (pre-shift) values for the moving columns ----
$snapD = @{}
$snapI = @{}
$snapJ = @{}
$snapK = @{}
$snapL = @{}
$snapM = @{}
$snapP = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, $colD).Value()
    $snapI[$r] = $ws.Cells.Item($r, $colI).Value()
    $snapJ[$r] = $ws.Cells.Item($r, $colJ).Value()
    $snapK[$r] = $ws.Cells.Item($r, $colK).Value()
    $snapL[$r] = $ws.Cells.Item($r, $colL).Value()
    $snapM[$r] = $ws.Cells.Item($r, $colM).Value()
    $snapP[$r] = $ws.Cells.Item($r, $colP).Value()
}

# --- 2) build row 1198 from the constant metadata of row 1197 -------------
$metaCols = 1,2,3,5,6,7,8,14,15,17,18   # A,B,C,E,F,G,H,N,O,Q,R
foreach ($c in $metaCols) {
    $ws.Cells.Item($newLastRow, $c).Value = $ws.Cells.Item($lastRow, $c).Value()
}
$ws.Cells.Item($newLastRow, $colD).NumberFormat = $ws.Cells.Item($lastRow, $colD).NumberFormat()

# --- 3) shift every record down by one row: new row r <- old row (r-1) ----
for ($r = $newLastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, $colD).Value = $snapD[$src]
    $ws.Cells.Item($r, $colI).Value = $snapI[$src]
    $ws.Cells.Item($r, $colJ).Value = $snapJ[$src]
    $ws.Cells.Item($r, $colK).Value = $snapK[$src]
    $ws.Cells.Item($r, $colL).Value = $snapL[$src]
    $ws.Cells.Item($r, $colM).Value = $snapM[$src]
    $ws.Cells.Item($r, $colP).Value = $snapP[$src]
}

# --- 4) write the brand-new record into row 1096 ---------------------------
# Quality/min/max/avg/price-per-kg stay the same as the former row 1096;
# only the date and the volume change.
$ws.Cells.Item($firstRow, $colD).Value = 45166
$ws.Cells.Item($firstRow, $colI).Value = $snapI[$firstRow]
$ws.Cells.Item($firstRow, $colJ).Value = 2000
$ws.Cells.Item($firstRow, $colK).Value = $snapK[$firstRow]
$ws.Cells.Item($firstRow, $colL).Value = $snapL[$firstRow]
$ws.Cells.Item($firstRow, $colM).Value = $snapM[$firstRow]
$ws.Cells.Item($firstRow, $colP).Value = $snapP[$firstRow]
